# Fill in match scores for matches that previously had empty H_score / A_score
# (the "13 October 2021" through "27 January 2022" match days), and widen the
# stored row span for the subsequent still-empty rows so it matches Excel's
# habit of recording the full A:E span once a sheet has been touched/saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scores")

# Row => (H_score, A_score)
$scores = @{
    22 = @(2,0)
    23 = @(2,3)
    24 = @(3,3)
    25 = @(1,2)
    26 = @(2,1)
    27 = @(1,0)
    28 = @(0,1)
    29 = @(0,0)
    30 = @(2,3)
    31 = @(1,1)
    32 = @(4,1)
    33 = @(1,2)
    34 = @(4,1)
    35 = @(0,4)
    36 = @(1,2)
    37 = @(2,0)
}

foreach ($row in $scores.Keys) {
    $pair = $scores[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# Update the view so it reflects where the user was working when saving
# (scrolled further down the sheet, with a new active cell selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E38").Select()
